# Apply updated crypto price/volume data to worksheet cells.
# D-column price strings that look numeric (e.g. "580.18") would be
# auto-converted to numeric values by COM Range.Value assignment, so we
# temporarily force a Text number format, assign the literal string, then
# clear the format again so the cell ends up with no explicit style --
# matching the original workbook (cells had no "s" attribute).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "63.744.17"
$ws.Range("E2").Value = "  -1.17%  "

Set-TextValue "D3" "2.639.05"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  -0.01%  "

Set-TextValue "D5" "580.18"
$ws.Range("E5").Value = "  -0.02%  "

Set-TextValue "D6" "155.48"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -4.64%  "

Set-TextValue "D9" "2.636.49"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("E10").Value = "  -4.15%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("E12").Value = "  -2.03%  "

$ws.Range("E13").Value = "  +0.92%  "

Set-TextValue "D14" "28.39"
$ws.Range("E14").Value = "  -1.03%  "

Set-TextValue "D15" "3.114.91"
$ws.Range("E15").Value = "  +0.23%  "

Set-TextValue "D16" "0.0000184"
$ws.Range("E16").Value = "  -1.56%  "

Set-TextValue "D17" "63.665.77"
$ws.Range("E17").Value = "  -0.97%  "

Set-TextValue "D18" "2.640.34"
$ws.Range("E18").Value = "  +0.23%  "

Set-TextValue "D19" "12.13"
$ws.Range("E19").Value = "  -1.22%  "

Set-TextValue "D20" "7.67"
$ws.Range("E20").Value = "  +3.53%  "

Set-TextValue "D21" "4.53"
$ws.Range("E21").Value = "  -3.18%  "

Set-TextValue "D22" "344.35"
$ws.Range("E22").Value = "  -0.71%  "

$ws.Range("E23").Value = "  +0.37%  "

Set-TextValue "D24" "68.02"
$ws.Range("E24").Value = "  -0.06%  "

Set-TextValue "D25" "1.89"
$ws.Range("E25").Value = "  +7.88%  "

$ws.Range("E26").Value = "  -3.88%  "

Set-TextValue "D27" "604.46"
$ws.Range("E27").Value = "  +3.74%  "

$ws.Range("E28").Value = "  -2.09%  "

$ws.Range("E29").Value = "  +1.88%  "

Set-TextValue "D30" "8.11"
$ws.Range("E30").Value = "  +2.28%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  -0.90%  "

Set-TextValue "D33" "2.07"
$ws.Range("E33").Value = "  -0.71%  "

$ws.Range("E34").Value = "  +0.83%  "

Set-TextValue "D35" "6.56"
$ws.Range("E35").Value = "  -1.84%  "

$ws.Range("E36").Value = "  +2.18%  "

Set-TextValue "D37" "0.403"
$ws.Range("E37").Value = "  -2.57%  "

$ws.Range("E38").Value = "  +0.01%  "

Set-TextValue "D39" "19.68"
$ws.Range("E39").Value = "  -2.15%  "

$ws.Range("E40").Value = "  -2.24%  "

Set-TextValue "D41" "151.00"
$ws.Range("E41").Value = "  -2.55%  "

$ws.Range("E42").Value = "  -0.02%  "

Set-TextValue "D43" "2.55"
$ws.Range("E43").Value = "  +3.26%  "

$ws.Range("E44").Value = "  -0.58%  "

Set-TextValue "D45" "161.41"
$ws.Range("E45").Value = "  +2.03%  "

Set-TextValue "D46" "24.26"
$ws.Range("E46").Value = "  +4.26%  "

Set-TextValue "D47" "3.90"
$ws.Range("E47").Value = "  -2.56%  "

$ws.Range("E48").Value = "  -2.68%  "

Set-TextValue "D49" "0.633"
$ws.Range("E49").Value = "  -0.64%  "

Set-TextValue "D50" "0.0998"
$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("E51").Value = "  -1.59%  "
